# Insert a new data row at row 295 (pushes the existing rows 295-358 down
# to 296-359, preserving all of their values/formatting) and populate the
# newly inserted row with a new "Pepino ensalada" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 295..358 down by one to make room for the new record.
$ws.Rows.Item(295).EntireRow.Insert()

# Fill in the new row 295 with the new observation's data.
$ws.Cells.Item(295, 1).Value  = 3
$ws.Cells.Item(295, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(295, 3).Value  = "Coquimbo"
$ws.Cells.Item(295, 4).Value  = 44711
$ws.Cells.Item(295, 5).Value  = 5
$ws.Cells.Item(295, 6).Value  = 100112043
$ws.Cells.Item(295, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(295, 8).Value  = "Sin especificar"
$ws.Cells.Item(295, 9).Value  = "Primera"
$ws.Cells.Item(295, 10).Value = 78
$ws.Cells.Item(295, 11).Value = 20000
$ws.Cells.Item(295, 12).Value = 21000
$ws.Cells.Item(295, 13).Value = 20487
$ws.Cells.Item(295, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(295, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(295, 16).Value = 293
$ws.Cells.Item(295, 17).Value = 70
$ws.Cells.Item(295, 18).Value = "Hortaliza"
